# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings (preserve as text, not auto-converted numbers)
# (applied per-cell in a loop since multi-area Range.NumberFormat assignment only affects the first area)
foreach ($r in @(5,6,7,8,9,10,11,12,15,17,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,43,47,48,49,50,51)) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "26.059.30"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.645.41"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "218.05"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "0.5189"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2616"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "0.06284"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "20.21"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "0.07660"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "4.562"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "1.616.85"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "1.873.12"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "0.5558"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "0.0₅8087"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "64.94"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "26.013.28"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "4.587"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "192.72"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "10.40"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").Value = "5.901"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "144.42"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").Value = "0.1175"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").Value = "7.165"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "15.80"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "0.05356"
$ws.Range("E30").Value = "  -5.38%  "
$ws.Range("D31").Value = "1.265"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "3.442"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "3.316"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "1.549"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").Value = "2.416"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "2.783"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "0.9373"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "0.5559"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").Value = "0.01569"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "1.004"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.746"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").Value = "1.025.38"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("D43").Value = "0.8237"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").Value = "1.782.94"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +6.57%  "
$ws.Range("D47").Value = "57.04"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "0.9987"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "0.4312"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "7.891"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "0.05096"
$ws.Range("E51").Value = "  -3.94%  "
